# Added error state propagation
# - general!B19 (tau_r) and general!B20 (tau_a): replace formula "=$B$5/2" with a
#   literal value of 55 (the downstream shared-formula cells in column E follow
#   automatically).
# - general!B22 (errorPropTestEnable): flip the flag from 0 to 1.
# - errorInjection sheet: populate the previously-placeholder gyro-bias-rate
#   (del_gy/del_gz) and accelerometer-bias (del_acclx/del_accly/del_acclz)
#   rows with real injected values + units, bump del_h, and tag del_gx with
#   its unit string too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "general"
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")

# tau_r: was "=$B$5/2" -> now a fixed literal
$general.Range("B19").Value = 55

# tau_a: was "=$B$5/2" -> now a fixed literal
$general.Range("B20").Value = 55

# errorPropTestEnable: 0 -> 1 (turn the test on)
$general.Range("B22").Value = 1

# Move the visible selection the way the author left it
$general.Range("C27").Select()

# ---------------------------------------------------------------------------
# Sheet "errorInjection"
# ---------------------------------------------------------------------------
$errInj = $wb.Worksheets.Item("errorInjection")

# New shared strings get appended to the string table in first-reference
# order, so touch "mg" (del_accl*) before "m/s3"/"m/s4" (del_gy/del_gz) to
# land on the same table layout the original author ended up with.

# del_acclx
$errInj.Range("B13").Value = 0.00327
$errInj.Range("C13").Value = "mg"

# del_accly
$errInj.Range("B14").Value = 0.00327
$errInj.Range("C14").Value = "mg"

# del_acclz
$errInj.Range("B15").Value = 0.00327
$errInj.Range("C15").Value = "mg"

# del_gx: keep the (now-correct) value, just add its units column
$errInj.Range("B9").Value = 0.00006
$errInj.Range("B9").NumberFormat = "0.00E+00"
$errInj.Range("C9").Value = "m/s2"

# del_gy
$errInj.Range("B10").Value = 0.00002
$errInj.Range("B10").NumberFormat = "0.00E+00"
$errInj.Range("C10").Value = "m/s3"

# del_gz
$errInj.Range("B11").Value = 0.00003
$errInj.Range("B11").NumberFormat = "0.00E+00"
$errInj.Range("C11").Value = "m/s4"

# del_h
$errInj.Range("B12").Value = 50

$errInj.Range("D7").Select()
